# Update microstate list: remove the replicate microstate "SM21_micro016"
# (v1.4.1) from the Sheet1 table, pulling the rows below it up by one and
# dropping the now-duplicated trailing row together with its 2D depiction
# picture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row that holds the replicate microstate ID.
$targetId = "SM21_micro016"
$foundRow = -1
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 3; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 2).Value2 -eq $targetId) {
        $foundRow = $r
        break
    }
}

if ($foundRow -gt 0) {
    # Shift the microstate ID / SMILES pairs below the replicate up by one
    # row. Each row keeps its own pre-existing formatting/style; only the
    # B (microstate ID) and C (SMILES) values move.
    for ($r = $foundRow; $r -lt $lastRow; $r++) {
        $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r + 1, 2).Value2
        $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r + 1, 3).Value2
    }

    # Drop the now-duplicated last row and its associated 2D depiction
    # picture (the last picture anchored to the sheet).
    $ws.Rows.Item($lastRow).Delete()
    $ws.Shapes.Item($ws.Shapes.Count).Delete()
}
